# Applies the F-column (想去人数 / "want to go" count) updates described in the
# commit's xlsx diff. Values come straight from the diff hunks, grouped by sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 33
$ws.Range("F3").Value = 8996
$ws.Range("F4").Value = 2699
$ws.Range("F6").Value = 323
$ws.Range("F7").Value = 837
$ws.Range("F8").Value = 716
$ws.Range("F9").Value = 131
$ws.Range("F10").Value = 80
$ws.Range("F12").Value = 909
$ws.Range("F13").Value = 3899
$ws.Range("F14").Value = 299
$ws.Range("F15").Value = 179
$ws.Range("F16").Value = 808
$ws.Range("F18").Value = 56
$ws.Range("F19").Value = 503
$ws.Range("F20").Value = 8
$ws.Range("F22").Value = 1409
$ws.Range("F23").Value = 1365
$ws.Range("F24").Value = 488
$ws.Range("F26").Value = 156
$ws.Range("F27").Value = 176
$ws.Range("F28").Value = 374
$ws.Range("F29").Value = 73
$ws.Range("F30").Value = 1023
$ws.Range("F33").Value = 735
$ws.Range("F36").Value = 101
$ws.Range("F38").Value = 31
$ws.Range("F40").Value = 226
$ws.Range("F41").Value = 199
$ws.Range("F42").Value = 382
$ws.Range("F43").Value = 28
$ws.Range("F44").Value = 30

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 5
$ws.Range("F6").Value = 1

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 33
$ws.Range("F4").Value = 323
$ws.Range("F5").Value = 837
$ws.Range("F6").Value = 716
$ws.Range("F7").Value = 131
$ws.Range("F8").Value = 80
$ws.Range("F10").Value = 909
$ws.Range("F12").Value = 3899
$ws.Range("F13").Value = 299
$ws.Range("F14").Value = 179
$ws.Range("F16").Value = 5
$ws.Range("F17").Value = 808
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = 56
$ws.Range("F22").Value = 503
$ws.Range("F23").Value = 8
$ws.Range("F26").Value = 1409
$ws.Range("F27").Value = 1365
$ws.Range("F28").Value = 488
$ws.Range("F30").Value = 156
$ws.Range("F31").Value = 176
$ws.Range("F33").Value = 374
$ws.Range("F34").Value = 73
$ws.Range("F35").Value = 1023
$ws.Range("F37").Value = 735
$ws.Range("F40").Value = 101
$ws.Range("F42").Value = 31
$ws.Range("F44").Value = 199
$ws.Range("F45").Value = 382
$ws.Range("F46").Value = 28
$ws.Range("F47").Value = 30
